$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-process dimension metadata for "municipio-nombre" (column C) and
# "personalidad-juridica" (column E) with the newly curated dimensions.

# municipio-nombre: now modeled as a proper "dim" (sdmx-dimension:refArea)
# whose value type is a municipality URI, instead of a measure of type int.
$ws.Range("C2").Value = "sdmx-dimension:refArea"
$ws.Range("C3").Value = "dim"
$ws.Range("C4").Value = "URI-Municipio"

# personalidad-juridica: now modeled as a "medida" (iaest-measure) of type
# xsd:int, instead of a dimension concept, and no longer needs a mapping file.
$ws.Range("E2").Value = "iaest-measure:personalidad-juridica"
$ws.Range("E3").Value = "medida"
$ws.Range("E4").Value = "xsd:int"
$ws.Range("E5").Clear()
